# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46057

$ws.Range("B2").Value = 8.029999999999999
$ws.Range("C2").Value = 1.45
$ws.Range("D2").Value = 0.05
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = -0.05
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0.02
$ws.Range("J2").Value = 2.58
$ws.Range("K2").Value = 10.73
$ws.Range("L2").Value = 9.949999999999999
$ws.Range("M2").Value = 5.94
$ws.Range("N2").Value = 3.78
$ws.Range("O2").Value = 3.78
$ws.Range("P2").Value = 3.78
$ws.Range("Q2").Value = 3.78
$ws.Range("R2").Value = 5.34
$ws.Range("S2").Value = 10.55
$ws.Range("T2").Value = 13.22
$ws.Range("U2").Value = 33.07
$ws.Range("V2").Value = 38.36
$ws.Range("W2").Value = 24.19
$ws.Range("X2").Value = 11.31
$ws.Range("Y2").Value = 6.99
$ws.Range("Z2").Value = 8.199999999999999

$ws.Range("AB2").Value = 20.21
$ws.Range("AD2").Value = 31.28
$ws.Range("AF2").Value = 23.14
